$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (N_Calib_1=20, N_Calib_2=40)
$ws.Range("C2").Value = 0.2322900381552812
$ws.Range("D2").Value = 0.8177053223704471

# Row 3 (N_Calib_1=20, N_Calib_2=60)
$ws.Range("C3").Value = 0.2721584051040664
$ws.Range("D3").Value = 0.7871462031477796

# Row 4 (N_Calib_1=20, N_Calib_2=100)
$ws.Range("C4").Value = -1.285616366838133
$ws.Range("D4").Value = 0.2072647646266859

# Row 5 (N_Calib_1=20, N_Calib_2=200)
$ws.Range("C5").Value = -1.886956868901056
$ws.Range("D5").Value = 0.06772870268374653
$ws.Range("G5").Value = "No"

# Row 6 (N_Calib_1=40, N_Calib_2=60)
$ws.Range("C6").Value = 0.04514201149405924
$ws.Range("D6").Value = 0.9642582122579351

# Row 7 (N_Calib_1=40, N_Calib_2=100)
$ws.Range("C7").Value = -1.277224288730588
$ws.Range("D7").Value = 0.2101765164717231

# Row 8 (N_Calib_1=40, N_Calib_2=200)
$ws.Range("C8").Value = -2.090731912133612
$ws.Range("D8").Value = 0.04409976294715556

# Row 9 (N_Calib_1=60, N_Calib_2=100)
$ws.Range("C9").Value = -1.309198466453444
$ws.Range("D9").Value = 0.1992463668175881

# Row 10 (N_Calib_1=60, N_Calib_2=200)
$ws.Range("C10").Value = -1.933660733967795
$ws.Range("D10").Value = 0.06151523306906181
$ws.Range("G10").Value = "No"

# Row 11 (N_Calib_1=100, N_Calib_2=200)
$ws.Range("C11").Value = -1.57740218652866
$ws.Range("D11").Value = 0.1239630975042803
$ws.Range("G11").Value = "No"
